$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 660.6429000000001
$ws.Range("J19").Value = 667.9474
$ws.Range("L19").Value = 667.9474
$ws.Range("N19").Value = -1017.9474
# Row 43
$ws.Range("H43").Value = 970.375
$ws.Range("I43").Value = 1095.25
$ws.Range("J43").Value = 845.5
$ws.Range("K43").Value = 1095.25
$ws.Range("L43").Value = 845.5
$ws.Range("M43").Value = -1026.25
$ws.Range("N43").Value = -983.5
# Row 137
$ws.Range("H137").Value = 1131.0952
$ws.Range("I137").Value = 870.0909
$ws.Range("J137").Value = 1418.2
$ws.Range("K137").Value = 2610.2727
$ws.Range("L137").Value = 4254.6
$ws.Range("M137").Value = -60.27269999999999
$ws.Range("N137").Value = -9354.6
# Row 138
$ws.Range("H138").Value = 9618235
$ws.Range("I138").Value = 4227.3335
$ws.Range("J138").Value = 12502437
$ws.Range("K138").Value = 12682.0005
$ws.Range("L138").Value = 37507311
$ws.Range("M138").Value = -7542.000499999998
$ws.Range("N138").Value = -37517591
# Row 141
$ws.Range("H141").Value = 2652.9167
$ws.Range("I141").Value = 2183.5
$ws.Range("K141").Value = 6550.5
$ws.Range("M141").Value = -1370.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 4
$ws.Range("H4").Value = 338.8
$ws.Range("I4").Value = 338.8
$ws.Range("K4").Value = 338.8
$ws.Range("M4").Value = -222.8
# Row 6
$ws.Range("H6").Value = 7024000
$ws.Range("I6").Value = 1280000
$ws.Range("J6").Value = 30000000
$ws.Range("K6").Value = 1280000
$ws.Range("L6").Value = 30000000
$ws.Range("M6").Value = -1279827
$ws.Range("N6").Value = -30000346
# Row 9
$ws.Range("H9").Value = 16250
$ws.Range("I9").Value = 50000
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 50000
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = -49830
$ws.Range("N9").Value = -5340
# Row 20
$ws.Range("H20").Value = 16250
$ws.Range("I20").Value = 50000
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 50000
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -49730
$ws.Range("N20").Value = -5540
# Row 23
$ws.Range("H23").Value = 21900
$ws.Range("J23").Value = 21900
$ws.Range("L23").Value = 21900
$ws.Range("N23").Value = -22418
# Row 56
$ws.Range("H56").Value = 10600
$ws.Range("J56").Value = 10600
$ws.Range("L56").Value = 10600
$ws.Range("N56").Value = -12084
# Row 74
$ws.Range("H74").Value = 639.2432
$ws.Range("I74").Value = 473.33334
$ws.Range("J74").Value = 945.53845
$ws.Range("K74").Value = 473.33334
$ws.Range("L74").Value = 945.53845
$ws.Range("M74").Value = 400.66666
$ws.Range("N74").Value = -2693.53845
# Row 77
$ws.Range("H77").Value = 639.2432
$ws.Range("I77").Value = 473.33334
$ws.Range("J77").Value = 945.53845
$ws.Range("K77").Value = 2366.6667
$ws.Range("L77").Value = 4727.69225
$ws.Range("M77").Value = 2001.3333
$ws.Range("N77").Value = -13463.69225
# Row 92
$ws.Range("H92").Value = 186516.67
$ws.Range("J92").Value = 186516.67
$ws.Range("L92").Value = 186516.67
$ws.Range("N92").Value = -191508.67
# Row 132
$ws.Range("H132").Value = 2606.8235
$ws.Range("I132").Value = 2224.05
$ws.Range("J132").Value = 3153.6428
$ws.Range("K132").Value = 6672.150000000001
$ws.Range("L132").Value = 9460.928400000001
$ws.Range("M132").Value = -4142.150000000001
$ws.Range("N132").Value = -14520.9284

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = 0
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
# Row 105
$ws.Range("H105").Value = 2993
$ws.Range("I105").Value = 2799.5652
$ws.Range("J105").Value = 3628.5715
$ws.Range("K105").Value = 2799.5652
$ws.Range("L105").Value = 3628.5715
$ws.Range("M105").Value = -1052.5652
$ws.Range("N105").Value = -7122.5715
# Row 107
$ws.Range("H107").Value = 2092.4443
$ws.Range("I107").Value = 2041.5
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 2041.5
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -121.5
$ws.Range("N107").Value = -6340

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1979.641
$ws.Range("I31").Value = 1328.375
$ws.Range("J31").Value = 2432.6956
$ws.Range("K31").Value = 1328.375
$ws.Range("L31").Value = 2432.6956
$ws.Range("M31").Value = -1033.375
$ws.Range("N31").Value = -3022.6956
# Row 34
$ws.Range("H34").Value = 1979.641
$ws.Range("I34").Value = 1328.375
$ws.Range("J34").Value = 2432.6956
$ws.Range("K34").Value = 1328.375
$ws.Range("L34").Value = 2432.6956
$ws.Range("M34").Value = -1126.375
$ws.Range("N34").Value = -2836.6956
# Row 92
$ws.Range("H92").Value = 30601
$ws.Range("J92").Value = 30601
$ws.Range("L92").Value = 30601
$ws.Range("N92").Value = -35593
# Row 99
$ws.Range("H99").Value = 8932285
$ws.Range("J99").Value = 2749.5
$ws.Range("L99").Value = 2749.5
$ws.Range("N99").Value = -5745.5
# Row 126
$ws.Range("H126").Value = 8932285
$ws.Range("J126").Value = 2749.5
$ws.Range("L126").Value = 8248.5
$ws.Range("N126").Value = -13188.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 92.666664
$ws.Range("I33").Value = 104.85714
$ws.Range("K33").Value = 629.14284
$ws.Range("M33").Value = -346.14284
# Row 38
$ws.Range("H38").Value = 141.20833
$ws.Range("J38").Value = 107.69231
$ws.Range("L38").Value = 323.07693
$ws.Range("N38").Value = -1017.07693
# Row 133
$ws.Range("H133").Value = 8166.6665
$ws.Range("I133").Value = 5900
$ws.Range("K133").Value = 17700
$ws.Range("M133").Value = -12640

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 26500
$ws.Range("I10").Value = 26500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 26500
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -26331
# Row 92
$ws.Range("H92").Value = 10125.5
$ws.Range("J92").Value = 10125.5
$ws.Range("L92").Value = 10125.5
$ws.Range("N92").Value = -13869.5
# Row 132
$ws.Range("H132").Value = 2709.6487
$ws.Range("I132").Value = 2104.7585
$ws.Range("J132").Value = 4902.375
$ws.Range("K132").Value = 6314.2755
$ws.Range("L132").Value = 14707.125
$ws.Range("M132").Value = -3784.2755
$ws.Range("N132").Value = -19767.125
# Row 138
$ws.Range("H138").Value = 62500
$ws.Range("J138").Value = 62500
$ws.Range("L138").Value = 62500
$ws.Range("N138").Value = -72780
# Row 139
$ws.Range("H139").Value = 59679
$ws.Range("J139").Value = 59679
$ws.Range("L139").Value = 59679
$ws.Range("N139").Value = -69959

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2925.7097
$ws.Range("I7").Value = 2079.7
$ws.Range("J7").Value = 3328.5715
$ws.Range("K7").Value = 2079.7
$ws.Range("L7").Value = 3328.5715
$ws.Range("M7").Value = -1967.7
$ws.Range("N7").Value = -3552.5715
# Row 12
$ws.Range("H12").Value = 2000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830
# Row 40
$ws.Range("H40").Value = 3128.9048
$ws.Range("I40").Value = 2415.2856
$ws.Range("K40").Value = 2415.2856
$ws.Range("M40").Value = -2279.2856
# Row 94
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31352
# Row 122
$ws.Range("H122").Value = 3708.0417
$ws.Range("I122").Value = 3284.7144
$ws.Range("J122").Value = 3882.353
$ws.Range("K122").Value = 9854.143199999999
$ws.Range("L122").Value = 11647.059
$ws.Range("M122").Value = -7404.143199999999
$ws.Range("N122").Value = -16547.059
# Row 126
$ws.Range("H126").Value = 2925.7097
$ws.Range("I126").Value = 2079.7
$ws.Range("J126").Value = 3328.5715
$ws.Range("K126").Value = 6239.099999999999
$ws.Range("L126").Value = 9985.7145
$ws.Range("M126").Value = -3769.099999999999
$ws.Range("N126").Value = -14925.7145

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 13515158
$ws.Range("I132").Value = 16130315
$ws.Range("K132").Value = 48390945
$ws.Range("M132").Value = -48388415
